# Update computed result values in the "pl_mw" results sheet for the
# "Case with 380 kV" run. Only the numeric values in columns B, C, D, E, F,
# K, L for data rows 2-25 (A = 0..23) change; everything else is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number, then new values for columns B, C, D, E, F, K, L (in that order).
$data = @(
    @(2,1.165568306009135,0.03342341788579262,0.3090660498763214,0.08600710581192672,4.618477121752392,0.6763136704454951,0.2225196855929852),
    @(3,1.145630774578677,0.03160631916440337,0.2967283011830375,0.08467592048555517,4.405935134467342,0.6589400703129229,0.2168610581896928),
    @(4,1.134448267733831,0.03046620734694017,0.2890757357962741,0.08388671296325256,4.275682007018474,0.6490472067143571,0.2135342625221028),
    @(5,1.130157057674609,0.02999531420353563,0.2859374313796224,0.08357216748358098,4.222660805723109,0.6452098223986411,0.2122155782463437),
    @(6,1.129460539229825,0.02991673854396737,0.2854151087297083,0.08352036362286341,4.213860048901353,0.644584325718327,0.2119988440226308),
    @(7,1.134389319739768,0.03045988237587594,0.2890334923379498,0.08388244231321096,4.274966715616387,0.6489946697318487,0.2135163285660582),
    @(8,1.158473629041367,0.03280185618523035,0.3048277452171249,0.08554226210956983,4.545138293844104,0.6701620266318002,0.2205378824118043),
    @(9,1.214139903148236,0.03720773664953114,0.3352084593897473,0.08902139005115473,5.077144509861711,0.7178550510727177,0.2354845402658015),
    @(10,1.260234718278099,0.04034064433398754,0.3571981080036437,0.09171581356577363,5.469705866913358,0.7567237436144296,0.2471937039239975),
    @(11,1.282345237154004,0.04174540326101095,0.3671363910969774,0.09297198724315336,5.648740027691247,0.7752508165879419,0.2526808745670763),
    @(12,1.290882992621363,0.0422745836958569,0.3708909178180591,0.09345207302565939,5.716606895914424,0.7823891224853696,0.2547819814463139),
    @(13,1.289036883633571,0.04216073583783242,0.3700827013003334,0.09334848195982204,5.701987335881938,0.7808462972135715,0.2543284348853661),
    @(14,1.28304433229286,0.04178899384586998,0.3674454537594727,0.09301139589445029,5.654322033977508,0.7758356287605181,0.2528532674357393),
    @(15,1.279395231182804,0.04156093514963999,0.3658289197511522,0.09280549453926668,5.625134991443588,0.7727824312117662,0.2519527147599518),
    @(16,1.258812786652101,0.04024844323768662,0.3565473498900076,0.09163433434171964,5.458015160744424,0.7555300479887137,0.2468383481200789),
    @(17,1.246479115466229,0.03943815943857487,0.3508370993792482,0.09092368210832191,5.355612851994181,0.7451634744682281,0.2437420841302895),
    @(18,1.239492532319531,0.03897017399987135,0.3475465768216708,0.09051780025866663,5.296756376690809,0.739280422529248,0.2419763041687588),
    @(19,1.237145420554526,0.03881138542246276,0.3464313946495423,0.09038086739766982,5.276835756767724,0.7373021532995665,0.2413810322733241),
    @(20,1.247780932725789,0.03952461426086984,0.351445597251967,0.09099903544245791,5.366509305552313,0.7462587765672026,0.2440701219826167),
    @(21,1.284800005839003,0.04189825737320518,0.3682203148378846,0.09311028667217869,5.668320535046007,0.7773040513268654,0.2532859281766804),
    @(22,1.309956138354266,0.04343349468108215,0.3791319824927086,0.09451576762990754,5.865984989641731,0.7983084318169631,0.2594444852621507),
    @(23,1.296441541160561,0.04261552578483219,0.3733127902926014,0.09376328201951978,5.7604482767847,0.7870323116867155,0.2561451029434068),
    @(24,1.247192056615745,0.03948553472171312,0.3511705191425847,0.09096495984909225,5.361582969509016,0.7457633512517248,0.2439217715279369),
    @(25,1.198172025981307,0.03603481746206683,0.3270499636984141,0.08805602403614188,4.932950692271675,0.7042847942202002,0.2313139945974285)
)

$columns = @("B", "C", "D", "E", "F", "K", "L")

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
